$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update timestamp in A2 ---
$ws.Range("A2").Value = "Laatst bijgewerkt: 2025-09-06 22:36:08"

# --- Update data rows 3-12 (columns D,E,F,G,H,I,J,K,L) ---
# Row 3
$ws.Range("D3").Value = "meer dan 10.5"
$ws.Range("E3").Value = "toto"
$ws.Range("F3").Value = 2.5
$ws.Range("G3").Value = "minder dan 10.5"
$ws.Range("H3").Value = "starcasino"
$ws.Range("I3").Value = 1.76
$ws.Range("J3").Value = "1=62, 2=88"
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = "€4.88"
$ws.Range("K3").Style = "Normal"
$ws.Range("L3").Value = 3.18

# Row 4
$ws.Range("D4").Value = "meer dan 9.5"
$ws.Range("E4").Value = "toto"
$ws.Range("F4").Value = 3.5
$ws.Range("G4").Value = "minder dan 9.5"
$ws.Range("H4").Value = "kambi"
$ws.Range("I4").Value = 1.44
$ws.Range("J4").Value = "1=44, 2=106"
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "€2.64"
$ws.Range("K4").Style = "Normal"
$ws.Range("L4").Value = 1.98

# Row 5
$ws.Range("D5").Value = "meer dan 8.5"
$ws.Range("E5").Value = "toto"
$ws.Range("F5").Value = 2.45
$ws.Range("G5").Value = "minder dan 8.5"
$ws.Range("H5").Value = "kambi"
$ws.Range("I5").Value = 1.74
$ws.Range("J5").Value = "1=62, 2=88"
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = "€1.9"
$ws.Range("K5").Style = "Normal"
$ws.Range("L5").Value = 1.71

# Row 6
$ws.Range("D6").Value = "meer dan 8.5"
$ws.Range("E6").Value = "toto"
$ws.Range("F6").Value = 2.45
$ws.Range("G6").Value = "minder dan 8.5"
$ws.Range("H6").Value = "jacks"
$ws.Range("I6").Value = 1.74
$ws.Range("J6").Value = "1=62, 2=88"
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = "€1.9"
$ws.Range("K6").Style = "Normal"
$ws.Range("L6").Value = 1.71

# Row 7
$ws.Range("D7").Value = "meer dan 8.5"
$ws.Range("E7").Value = "toto"
$ws.Range("F7").Value = 2.45
$ws.Range("G7").Value = "minder dan 8.5"
$ws.Range("H7").Value = "betmgm"
$ws.Range("I7").Value = 1.74
$ws.Range("J7").Value = "1=62, 2=88"
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "€1.9"
$ws.Range("K7").Style = "Normal"
$ws.Range("L7").Value = 1.71

# Row 8
$ws.Range("D8").Value = "meer dan 7.5"
$ws.Range("E8").Value = "toto"
$ws.Range("F8").Value = 1.85
$ws.Range("G8").Value = "minder dan 7.5"
$ws.Range("H8").Value = "betmgm"
$ws.Range("I8").Value = 2.25
$ws.Range("J8").Value = "1=82, 2=68"
$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = "€1.7"
$ws.Range("K8").Style = "Normal"
$ws.Range("L8").Value = 1.5

# Row 9
$ws.Range("D9").Value = "meer dan 7.5"
$ws.Range("E9").Value = "toto"
$ws.Range("F9").Value = 1.85
$ws.Range("G9").Value = "minder dan 7.5"
$ws.Range("H9").Value = "jacks"
$ws.Range("I9").Value = 2.25
$ws.Range("J9").Value = "1=82, 2=68"
$ws.Range("K9").NumberFormat = "@"
$ws.Range("K9").Value = "€1.7"
$ws.Range("K9").Style = "Normal"
$ws.Range("L9").Value = 1.5

# Row 10
$ws.Range("D10").Value = "meer dan 7.5"
$ws.Range("E10").Value = "toto"
$ws.Range("F10").Value = 1.85
$ws.Range("G10").Value = "minder dan 7.5"
$ws.Range("H10").Value = "kambi"
$ws.Range("I10").Value = 2.25
$ws.Range("J10").Value = "1=82, 2=68"
$ws.Range("K10").NumberFormat = "@"
$ws.Range("K10").Value = "€1.7"
$ws.Range("K10").Style = "Normal"
$ws.Range("L10").Value = 1.5

# Row 11
$ws.Range("D11").Value = "meer dan 9.5"
$ws.Range("E11").Value = "toto"
$ws.Range("F11").Value = 3.5
$ws.Range("G11").Value = "minder dan 9.5"
$ws.Range("H11").Value = "betmgm"
$ws.Range("I11").Value = 1.43
$ws.Range("J11").Value = "1=44, 2=106"
$ws.Range("K11").NumberFormat = "@"
$ws.Range("K11").Value = "€1.58"
$ws.Range("K11").Style = "Normal"
$ws.Range("L11").Value = 1.5

# Row 12
$ws.Range("D12").Value = "meer dan 9.5"
$ws.Range("E12").Value = "toto"
$ws.Range("F12").Value = 3.5
$ws.Range("G12").Value = "minder dan 9.5"
$ws.Range("H12").Value = "jacks"
$ws.Range("I12").Value = 1.43
$ws.Range("J12").Value = "1=44, 2=106"
$ws.Range("K12").NumberFormat = "@"
$ws.Range("K12").Value = "€1.58"
$ws.Range("K12").Style = "Normal"
$ws.Range("L12").Value = 1.5

# --- Set M/N column text values (hyperlink display text) ---
$ws.Range("M3").Value = "https://sport.toto.nl/wedden/wedstrijd/8778584"
$ws.Range("N3").Value = "https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=12642394"
$ws.Range("M4").Value = "https://sport.toto.nl/wedden/wedstrijd/8706282"
$ws.Range("N4").Value = "https://www.unibet.nl/betting/sports/event/1023224945?coupon=single%7C3865617128%7C1.44%7Creplace"
$ws.Range("M5").Value = "https://sport.toto.nl/wedden/wedstrijd/8706282"
$ws.Range("N5").Value = "https://www.unibet.nl/betting/sports/event/1023224945?coupon=single%7C3865617143%7C1.74%7Creplace"
$ws.Range("M6").Value = "https://sport.toto.nl/wedden/wedstrijd/8706282"
$ws.Range("N6").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"
$ws.Range("M7").Value = "https://sport.toto.nl/wedden/wedstrijd/8706282"
$ws.Range("N7").Value = "https://www.betmgm.nl/betting/sports/event/1023224945?coupon=single%7C3865617143%7C1.74%7Creplace"
$ws.Range("M8").Value = "https://sport.toto.nl/wedden/wedstrijd/8706282"
$ws.Range("N8").Value = "https://www.betmgm.nl/betting/sports/event/1023224945?coupon=single%7C3865617122%7C2.25%7Creplace"
$ws.Range("M9").Value = "https://sport.toto.nl/wedden/wedstrijd/8706282"
$ws.Range("N9").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"
$ws.Range("M10").Value = "https://sport.toto.nl/wedden/wedstrijd/8706282"
$ws.Range("N10").Value = "https://www.unibet.nl/betting/sports/event/1023224945?coupon=single%7C3865617122%7C2.25%7Creplace"
$ws.Range("M11").Value = "https://sport.toto.nl/wedden/wedstrijd/8706282"
$ws.Range("N11").Value = "https://www.betmgm.nl/betting/sports/event/1023224945?coupon=single%7C3865617128%7C1.43%7Creplace"
$ws.Range("M12").Value = "https://sport.toto.nl/wedden/wedstrijd/8706282"
$ws.Range("N12").Value = "https://jacks.nl/sports/event/1023224945#event/1023224945"

# --- Rebuild hyperlinks (delete all, re-add with correct targets/fragments) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("M3"), "https://sport.toto.nl/wedden/wedstrijd/8778584", "", "", "")
$ws.Hyperlinks.Add($ws.Range("N3"), "https://starcasino.nl/prematch-bets?page=event&sportId=66&eventId=12642394", "", "", "")
$ws.Hyperlinks.Add($ws.Range("M4"), "https://sport.toto.nl/wedden/wedstrijd/8706282", "", "", "")
$ws.Hyperlinks.Add($ws.Range("N4"), "https://jacks.nl/sports/event/1023224945", "", "", "")
$ws.Hyperlinks.Add($ws.Range("M5"), "https://sport.toto.nl/wedden/wedstrijd/8706282", "", "", "")
$ws.Hyperlinks.Add($ws.Range("N5"), "https://www.unibet.nl/betting/sports/event/1023224945?coupon=single%7C3865617143%7C1.74%7Creplace", "", "", "")
$ws.Hyperlinks.Add($ws.Range("M6"), "https://sport.toto.nl/wedden/wedstrijd/8706282", "", "", "")
$ws.Hyperlinks.Add($ws.Range("N6"), "https://www.betmgm.nl/betting/sports/event/1023224945?coupon=single%7C3865617143%7C1.74%7Creplace", "event/1023224945", "", "")
$ws.Hyperlinks.Add($ws.Range("M7"), "https://sport.toto.nl/wedden/wedstrijd/8706282", "", "", "")
$ws.Hyperlinks.Add($ws.Range("N7"), "https://www.betmgm.nl/betting/sports/event/1023224945?coupon=single%7C3865617122%7C2.25%7Creplace", "", "", "")
$ws.Hyperlinks.Add($ws.Range("M8"), "https://sport.toto.nl/wedden/wedstrijd/8706282", "", "", "")
$ws.Hyperlinks.Add($ws.Range("N8"), "https://www.unibet.nl/betting/sports/event/1023224945?coupon=single%7C3865617122%7C2.25%7Creplace", "", "", "")
$ws.Hyperlinks.Add($ws.Range("M9"), "https://sport.toto.nl/wedden/wedstrijd/8706282", "", "", "")
$ws.Hyperlinks.Add($ws.Range("N9"), "https://jacks.nl/sports/event/1023224945", "event/1023224945", "", "")
$ws.Hyperlinks.Add($ws.Range("M10"), "https://sport.toto.nl/wedden/wedstrijd/8706282", "", "", "")
$ws.Hyperlinks.Add($ws.Range("N10"), "https://www.betmgm.nl/betting/sports/event/1023224945?coupon=single%7C3865617128%7C1.43%7Creplace", "", "", "")
$ws.Hyperlinks.Add($ws.Range("M11"), "https://sport.toto.nl/wedden/wedstrijd/8706282", "", "", "")
$ws.Hyperlinks.Add($ws.Range("N11"), "https://jacks.nl/sports/event/1023224945", "", "", "")
$ws.Hyperlinks.Add($ws.Range("M12"), "https://sport.toto.nl/wedden/wedstrijd/8706282", "", "", "")
$ws.Hyperlinks.Add($ws.Range("N12"), "https://www.unibet.nl/betting/sports/event/1023224945?coupon=single%7C3865617128%7C1.43%7Creplace", "event/1023224945", "", "")

# --- Restore Hyperlink style on M/N cells (Add() bumps to an auto style otherwise) ---
$ws.Range("M3").Style = "Hyperlink"
$ws.Range("N3").Style = "Hyperlink"
$ws.Range("M4").Style = "Hyperlink"
$ws.Range("N4").Style = "Hyperlink"
$ws.Range("M5").Style = "Hyperlink"
$ws.Range("N5").Style = "Hyperlink"
$ws.Range("M6").Style = "Hyperlink"
$ws.Range("N6").Style = "Hyperlink"
$ws.Range("M7").Style = "Hyperlink"
$ws.Range("N7").Style = "Hyperlink"
$ws.Range("M8").Style = "Hyperlink"
$ws.Range("N8").Style = "Hyperlink"
$ws.Range("M9").Style = "Hyperlink"
$ws.Range("N9").Style = "Hyperlink"
$ws.Range("M10").Style = "Hyperlink"
$ws.Range("N10").Style = "Hyperlink"
$ws.Range("M11").Style = "Hyperlink"
$ws.Range("N11").Style = "Hyperlink"
$ws.Range("M12").Style = "Hyperlink"
$ws.Range("N12").Style = "Hyperlink"

Write-Output "edit complete"